$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting (style) from the existing
# header cell G1 so it matches the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values (1) for the data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
